# NB13 + NB15 angepasst
# - Kosten!B2 value changes from 33850 to 53550 (NB13)
# - Rows 3-14 are removed, shrinking the data table to just the header row
#   and the single data row (A2:C2) (NB15)
# - The line chart's two series (category + value refs) are updated to
#   point at the single remaining data row instead of the old A2:A14 /
#   B2:B14 / C2:C14 ranges.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Kosten")

# Update the chart series so they reference only the remaining row (row 2)
# instead of the soon-to-be-deleted rows 3-14.
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart

$s1 = $chart.SeriesCollection(1)
$s1.Formula = "=SERIES(Kosten!`$B`$1,Kosten!`$A`$2,Kosten!`$B`$2,1)"

$s2 = $chart.SeriesCollection(2)
$s2.Formula = "=SERIES(Kosten!`$C`$1,Kosten!`$A`$2,Kosten!`$C`$2,2)"

# Remove rows 3 through 14 entirely (shifts everything below them up,
# shrinking the sheet's used range / dimension to A1:C2).
$ws.Range("A3:C14").EntireRow.Delete()

# Adjust the remaining data row's "Beste Werte" value.
$ws.Range("B2").Value = 53550
